# Generate Report for handoff
#
# Semantic summary of the change:
#  - The localization source file "7a1823cc-6862-4f4d-b6cf-062a8cc9746d.md" was
#    replaced by a newer revision "3d7a97e8-37f9-49f3-852b-8681dab3bbc0.md", with
#    fresh handoff (.xlf) artifacts and timestamps for both the zh-cn and de-de
#    targets.
#  - A second source file "7ac6053b-ae52-4416-badd-f3e7d8b60527.md" showed up
#    whose handoff transform failed, so a new row is added for it (status
#    "Handoff transform failed") ahead of the existing ".localization-config"
#    bookkeeping row on every sheet.

$wb = $excel.ActiveWorkbook

$file1New  = "3d7a97e8-37f9-49f3-852b-8681dab3bbc0.md"
$file2New  = "7ac6053b-ae52-4416-badd-f3e7d8b60527.md"
$config    = ".localization-config"

$xlfZhNew  = "3d7a97e8-37f9-49f3-852b-8681dab3bbc0.cc935fac8bcabfef3627a6a6cd1ab5208981a7e0.zh-cn.xlf"
$xlfDeNew  = "3d7a97e8-37f9-49f3-852b-8681dab3bbc0.cc935fac8bcabfef3627a6a6cd1ab5208981a7e0.de-de.xlf"

$dtZhNew   = "2016-01-25 11:11:59"
$dtDeNew   = "2016-01-25 11:12:10"
$epoch     = "0001-01-01 00:00:00"

$urlMdBase     = "https://github.com/OpenLocalizationTest/oltest/blob/a9b4b6d1dbdda8dba524349d6c752163e50ea66e/e2e/"
$urlConfig     = "https://github.com/OpenLocalizationTest/oltest/blob/a9b4b6d1dbdda8dba524349d6c752163e50ea66e/.localization-config"
$urlXlfZhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73a340da15e9a037b760cdc9a7d9cf8273c9a42b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/"
$urlXlfDeBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff34f83b2235309f9800d33e2e8036b282dd0150/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Cells.Item(2,1).Value = $file1New
$wsOverview.Cells.Item(2,2).Value = "Ready for handoff"
$wsOverview.Cells.Item(2,3).Value = "Ready for handoff"

$wsOverview.Cells.Item(3,1).Value = $file2New
$wsOverview.Cells.Item(3,2).Value = "Handoff transform failed"
$wsOverview.Cells.Item(3,3).Value = "Handoff transform failed"

$wsOverview.Cells.Item(4,1).Value = $config
$wsOverview.Cells.Item(4,2).Value = "Not to be localized"
$wsOverview.Cells.Item(4,3).Value = "Not to be localized"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($urlMdBase + $file1New), "", "", $file1New)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($urlMdBase + $file2New), "", "", $file2New)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $urlConfig, "", "", $config)

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Cells.Item(2,1).Value = $file1New
$wsZh.Cells.Item(2,2).Value = "Ready for handoff"
$wsZh.Cells.Item(2,3).Value = $xlfZhNew
$wsZh.Cells.Item(2,4).Value = $dtZhNew
$wsZh.Cells.Item(2,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(2,7).Value = $epoch
$wsZh.Cells.Item(2,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(2,8).Value = "Include"

$wsZh.Cells.Item(3,1).Value = $file2New
$wsZh.Cells.Item(3,2).Value = "Handoff transform failed"

$wsZh.Cells.Item(4,1).Value = $config
$wsZh.Cells.Item(4,2).Value = "Not to be localized"
$wsZh.Cells.Item(4,4).Value = $epoch
$wsZh.Cells.Item(4,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,7).Value = $epoch
$wsZh.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,8).Value = "Ignored"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($urlMdBase + $file1New), "", "", $file1New)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), ($urlXlfZhBase + $xlfZhNew), "", "", $xlfZhNew)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($urlMdBase + $file2New), "", "", $file2New)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $urlConfig, "", "", $config)

# ---------------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Cells.Item(2,1).Value = $file1New
$wsDe.Cells.Item(2,2).Value = "Ready for handoff"
$wsDe.Cells.Item(2,3).Value = $xlfDeNew
$wsDe.Cells.Item(2,4).Value = $dtDeNew
$wsDe.Cells.Item(2,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(2,7).Value = $epoch
$wsDe.Cells.Item(2,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(2,8).Value = "Include"

$wsDe.Cells.Item(3,1).Value = $file2New
$wsDe.Cells.Item(3,2).Value = "Handoff transform failed"

$wsDe.Cells.Item(4,1).Value = $config
$wsDe.Cells.Item(4,2).Value = "Not to be localized"
$wsDe.Cells.Item(4,4).Value = $epoch
$wsDe.Cells.Item(4,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,7).Value = $epoch
$wsDe.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,8).Value = "Ignored"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($urlMdBase + $file1New), "", "", $file1New)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), ($urlXlfDeBase + $xlfDeNew), "", "", $xlfDeNew)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($urlMdBase + $file2New), "", "", $file2New)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $urlConfig, "", "", $config)
